$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.478.26"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "2.456.63"
$ws.Range("E3").Value = "  +1.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "580.68"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.97%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.00"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.44%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.78%  "
$ws.Range("D9").Value = "2.452.65"
$ws.Range("E9").Value = "  +1.03%  "
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("E11").Value = "  +2.48%  "
$ws.Range("E12").Value = "  +1.21%  "
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.53"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.36%  "
$ws.Range("D16").Value = "2.897.92"
$ws.Range("E16").Value = "  +1.24%  "
$ws.Range("D17").Value = "62.148.55"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "2.454.96"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.89"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "328.95"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +1.51%  "
$ws.Range("E22").Value = "  +0.55%  "
$ws.Range("E23").Value = "  -2.70%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "65.78"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.30"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.36%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "589.47"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.79%  "
$ws.Range("D28").Value = "0.0₃0975"
$ws.Range("E28").Value = "  +0.76%  "
$ws.Range("D29").Value = "2.578.13"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("E31").Value = "  -0.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.08"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.71%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.91"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.44%  "
$ws.Range("E34").Value = "  +0.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.96"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.81%  "
$ws.Range("E36").Value = "  +0.17%  "
$ws.Range("E37").Value = "  -0.92%  "
$ws.Range("E38").Value = "  +1.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "154.31"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +4.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.35"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.49"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.03%  "
$ws.Range("E42").Value = "  +2.55%  "
$ws.Range("E43").Value = "  +0.38%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.50"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "143.57"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("E47").Value = "  -0.81%  "
$ws.Range("D48").Value = "0.0₆0264"
$ws.Range("E48").Value = "  +21.19%  "
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0528"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.96"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.16%  "
